$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 - Le, Jimmy - add Assignment 1 and Assignment 2 scores
$ws.Range("D19").Value = 34
$ws.Range("E19").Value = 19

# Row 20 - Liang, Jia Q. - update Assignment 1, add Assignment 2 and 3 scores
$ws.Range("D20").Value = 41
$ws.Range("E20").Value = 19
$ws.Range("F20").Value = 24

# Row 24 - Navarro, Ignacio - add Assignment 1, 2, 3 scores
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = 19
$ws.Range("F24").Value = 24

# Row 29 - Quijano, Jesse A. - add Assignment 3 score
$ws.Range("F29").Value = 24

# Row 39 - Valino, Joshua F. - add Assignment 1 score
$ws.Range("D39").Value = 41

# Update the active selection to match the last edited cell
$ws.Range("D39").Select()
